$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Overview" to "ExtremeFlows"
$ws.Name = "ExtremeFlows"

# Update row 8-10 labels (indentation increased from 3 to 7 leading spaces)
$ws.Range("A8").Value = "       Lowest consecutive 10-year flow"
$ws.Range("A9").Value = "       Lowest consecutive 4-year flow"
$ws.Range("A10").Value = "       Lowest consecutive 3-year flow."

# Fill in row 7 (previously blank D/E/F) - match the existing vertical-top alignment
# used by the rest of the row before typing the quote-prefixed empty text value
$ws.Range("F7").VerticalAlignment = -4160
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "'"

# Remove the now unused rows (bottom-up so row numbers don't shift unexpectedly)
$ws.Rows(17).Delete()
$ws.Rows(11).Delete()

[void]$ws.Range("F7").Select()
